# "excel file and Base test class is done"
# Rename Sheet1 -> ValidLogin and populate it with the UserName/Password/
# HomePageTitle sample data used by the actiTIME login test, then tidy up
# the sheet (bold header row, widened title column, printable page setup).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ValidLogin"

# Header row (bold)
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "HomePageTitle"
$ws.Range("A1:C1").Font.Bold = $true

# Data row
$ws.Range("A2").Value = "purnendukumar82"
$ws.Range("B2").Value = 8706
$ws.Range("C2").Value = "actiTIME-Enter Time-Track"

# Widen column C so the long title fits
$ws.Columns.Item(3).ColumnWidth = 24.42578125

# Page setup used for printing this sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor below the data, matching the saved selection
[void]$ws.Range("A3").Select()
